# Correct a handful of "waste" (L column) measurements on the "home"
# sheet. The dependent ratio/evaluation columns (M:P) are formulas and
# recalculate automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("home")
$ws.Activate()

$ws.Range("L9").Value  = 560
$ws.Range("L11").Value = 1076
$ws.Range("L12").Value = 560
$ws.Range("L14").Value = 560
$ws.Range("L17").Value = 560

# Match the authored selection state (cell C10 active on the home tab).
$ws.Range("C10").Select()

$wb.Save()
